# "error solve ifrs list"
# Fix the IFRS financial data table on the active sheet:
#  - Rows 2-6 (2015/12 .. 2019/12 columns of source data) get corrected
#    figures, and columns J (당기순이익(비지배)) and O (자본총계(비지배))
#    are removed for those rows.
#  - Rows 7-9 (2020(E)/2021(E)/2022(E) placeholder rows) had been filled
#    in with erroneous data; all of their figures (columns D through AJ)
#    are cleared, leaving only the row label columns (A, B, C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New correct values per row; keyed by 1-based column index.
# Column indices: D=4 E=5 F=6 G=7 H=8 I=9 J=10(removed) K=11 L=12 M=13
# N=14 O=15(removed) P=16 Q=17 R=18 S=19 T=20 U=21 V=22 W=23 X=24 Y=25
# Z=26 AA=27 AB=28 AC=29 AD=30 AE=31 AF=32 AG=33 AH=34 AI=35 AJ=36
$rowValues = @{
    2 = @{ 4 = 1861; 5 = 101; 6 = 101; 7 = 89; 8 = 89; 9 = 89; 11 = 893; 12 = 606; 13 = 286; 14 = 286; 16 = 23; 17 = 0; 18 = -44; 19 = 94; 20 = 44; 21 = -44; 22 = 444; 23 = 5.45; 24 = 4.78; 25 = 40.4; 26 = 10.96; 27 = 211.71; 28 = 487.11; 29 = 831; 30 = 5.21; 31 = 2688; 32 = 1.61; 33 = 0; 34 = 0; 35 = 0; 36 = 10869068 }
    3 = @{ 4 = 1703; 5 = 87; 6 = 87; 7 = 73; 8 = 52; 9 = 52; 11 = 938; 12 = 625; 13 = 312; 14 = 312; 16 = 23; 17 = 39; 18 = -34; 19 = 17; 20 = 33; 21 = 6; 22 = 448; 23 = 5.08; 24 = 3.07; 25 = 17.49; 26 = 5.72; 27 = 200.26; 28 = 751.98; 29 = 482; 30 = 9.09; 31 = 2931; 32 = 1.49; 33 = 0; 34 = 0; 35 = 0; 36 = 10869068 }
    4 = @{ 4 = 1974; 5 = 106; 6 = 106; 7 = 97; 8 = 62; 9 = 62; 11 = 1241; 12 = 724; 13 = 517; 14 = 517; 16 = 30; 17 = 7; 18 = -159; 19 = 172; 20 = 117; 21 = -109; 22 = 502; 23 = 5.35; 24 = 3.15; 25 = 14.97; 26 = 5.7; 27 = 139.97; 28 = 1185.85; 29 = 481; 30 = 7.2; 31 = 3840; 32 = 0.9; 33 = 0; 34 = 0; 35 = 0; 36 = 13687934 }
    5 = @{ 4 = 1916; 5 = -59; 6 = -59; 7 = -99; 8 = -90; 9 = -90; 11 = 1063; 12 = 637; 13 = 426; 14 = 426; 16 = 65; 17 = -68; 18 = -44; 19 = 30; 20 = 21; 21 = -89; 22 = 439; 23 = -3.06; 24 = -4.68; 25 = -19.01; 26 = -7.78; 27 = 149.37; 28 = 414.79; 29 = -613; 30 = -4.32; 31 = 2817; 32 = 0.9399999999999999; 33 = 0; 34 = 0; 35 = 0; 36 = 15350354 }
    6 = @{ 4 = 2196; 5 = 57; 6 = 57; 7 = 42; 8 = 30; 9 = 30; 11 = 1230; 12 = 750; 13 = 479; 14 = 479; 16 = 69; 17 = -39; 18 = 15; 19 = 8; 20 = 14; 21 = -53; 22 = 444; 23 = 2.59; 24 = 1.35; 25 = 6.53; 26 = 2.58; 27 = 156.47; 28 = 446.88; 29 = 185; 30 = 12.36; 31 = 3003; 32 = 0.76; 33 = 0; 34 = 0; 35 = 0; 36 = 16183909 }
}

foreach ($rowNum in $rowValues.Keys) {
    $cols = $rowValues[$rowNum]
    foreach ($colNum in $cols.Keys) {
        $ws.Cells.Item($rowNum, $colNum).Value = $cols[$colNum]
    }
}

# Columns J and O no longer apply to rows 2-6; remove their stray values.
$ws.Range("J2:J6").ClearContents()
$ws.Range("O2:O6").ClearContents()

# Rows 7-9 were mistakenly populated with data; clear everything except
# the leading label columns A (순번), B (연간), and C (기간 구분).
$ws.Range("D7:AJ9").ClearContents()
